$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.348.98"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.602.76"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.007"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.281"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08136"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.678"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.583"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001270"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "1.596.27"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06840"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.636"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("D24").Value = "23.351.79"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.918"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.506"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.457"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("D33").Value = "1.776.18"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07759"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9591"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02803"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.361"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2557"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.01%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08878"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.389"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7199"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6623"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.324"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.988"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08036"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("E51").Value = "  -4.91%  "
